$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select the entire second row (the CO number that has been used/consumed)
# then delete it, which shifts all rows below it up by one.
$ws.Rows("2:2").Select()
$ws.Rows("2:2").Delete()
